$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 44, shifting existing rows 44-69 down to 45-70
$ws.Rows.Item(44).Insert()

# Populate the newly inserted row 44 with the new weekly data entry
$ws.Range("A44").Value = 11
$ws.Range("B44").Value = "Vega Monumental Concepción"
$ws.Range("C44").Value = "Bíobío"
$ws.Range("D44").Value = 44806
$ws.Range("E44").Value = 8
$ws.Range("F44").Value = 100112013
$ws.Range("G44").Value = "Alcachofa"
$ws.Range("H44").Value = "Argentina(o)"
$ws.Range("I44").Value = "Primera"
$ws.Range("J44").Value = 220
$ws.Range("K44").Value = 10000
$ws.Range("L44").Value = 12000
$ws.Range("M44").Value = 10909
$ws.Range("N44").Value = "$/caja 50 unidades"
$ws.Range("O44").Value = "Provincia de Limarí"
$ws.Range("P44").Value = 218
$ws.Range("Q44").Value = 50
$ws.Range("R44").Value = "Hortaliza"
